$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "2022" column (S) -------------------------------------
# Reuse R3 / R4's existing cell formats (rather than Range.Style=, which
# round-trips through a different path and mints brand-new style records).
$ws.Range("R3").Copy()
$ws.Range("S3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 0.071025550219041236

# --- Unify column widths for A:C ----------------------------------------
$ws.Range("A1:C1").ColumnWidth = 32.6

# --- Update the active selection ----------------------------------------
$ws.Range("F14").Select()
